$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# studyDesignProcedures: insert two new columns (procedureName,
# procedureDescription) between the existing procedureType and
# procedureCode columns, and populate them with data.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("studyDesignProcedures")

# Insert two blank columns before column C. This shifts the previous
# C/D/E columns (procedureCode/procedureIsConditional/
# procedureIsConditionalReason) to E/F/G and lets the new columns inherit
# the neighbouring column formatting automatically.
$ws.Columns("C:D").Insert()

# Match the author's column widths for the two new columns.
$ws.Range("C1").ColumnWidth = 16.666666666666668
$ws.Range("D1").ColumnWidth = 20.666666666666668

# Stage the new header/body text off to the side so that writing the
# values doesn't disturb the cell formatting that Insert() already set up
# for the new C/D columns (direct .Value assignment on a freshly
# formatted cell can reset some format bits, so we fill via copy/paste of
# values only).
$ws.Range("J1").Value = "procedureName"
$ws.Range("J2").Value = "procedureDescription"
$ws.Range("J3").Value = "Test8"
$ws.Range("J4").Value = "Test9"
$ws.Range("J5").Value = "Test Eight"
$ws.Range("J6").Value = "Test Nine"

$ws.Range("J1").Copy()
$ws.Range("C1").PasteSpecial(-4163)
$ws.Range("J2").Copy()
$ws.Range("D1").PasteSpecial(-4163)
$ws.Range("J3").Copy()
$ws.Range("C2").PasteSpecial(-4163)
$ws.Range("J4").Copy()
$ws.Range("C3").PasteSpecial(-4163)
$ws.Range("J5").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("J6").Copy()
$ws.Range("D3").PasteSpecial(-4163)

# Remove the scratch values.
$ws.Range("J1:J6").Clear()

# Make this the active sheet/cell, matching the author's final selection.
$ws.Activate()
$ws.Range("D6").Select()
